# Commit: "Test cases with the Bioblank and Diagnosis filters"
#
# The four Cypher queries stored in column B (rows 2-5, one per tab:
# CasesTab, SamplesTab, FilesTab, StudyFilesTab) each get an
# "Order by ... LIMIT 100" clause appended to their RETURN statement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - CasesTab query (cell B2)
$b2 = $ws.Cells.Item(2, 2).Value()
$ws.Cells.Item(2, 2).Value = $b2 + "`n       Order by c.case_id LIMIT 100"

# Row 3 - SamplesTab query (cell B3)
$b3 = $ws.Cells.Item(3, 2).Value()
$ws.Cells.Item(3, 2).Value = $b3 + "`n        Order by samp.sample_id LIMIT 100"

# Row 4 - FilesTab query (cell B4)
$b4 = $ws.Cells.Item(4, 2).Value()
$ws.Cells.Item(4, 2).Value = $b4 + "`n        Order By f.file_name LIMIT 100"

# Row 5 - StudyFilesTab query (cell B5)
$b5 = $ws.Cells.Item(5, 2).Value()
$ws.Cells.Item(5, 2).Value = $b5 + "`n    order by 'File Name' asc`n  limit 100"

# Appending text to the wrap-text query cells re-triggers row autofit in this
# host; pin rows 3-5 back to their fixed 100pt custom height (row 2 has no
# customHeight flag and is left alone).
$ws.Rows.Item(3).RowHeight = 100
$ws.Rows.Item(4).RowHeight = 100
$ws.Rows.Item(5).RowHeight = 100

# Sheet view / window tweaks (zoom level + selection moved to B6)
$ws.Application.ActiveWindow.Zoom = 80
$ws.Range("B6").Select()
